# Apply the "fix(gui) step 1 and 2" changes to the BUZON price list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1 by one day (45308 -> 45309, i.e. 2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2: update the prices for the B-009S and B-010 items
$ws.Range("D29").Value = 19600
$ws.Range("D30").Value = 21660
